# Refresh "ESTADO DE CUENTA" worker mora-period table (B15:J23):
# the previous account-statement periods are replaced by the new set -
# the "Periodo Mora" column is now listed most-recent-period-first, and
# the "Valor Mora" amount that used to sit on the last (2106) period now
# belongs to the first (2106) period, while the first row's old amount
# moves to the last (2011) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2106"
$ws.Range("E17").Value = "2105"
$ws.Range("E18").Value = "2104"
$ws.Range("E19").Value = "2103"
$ws.Range("E20").Value = "2102"
$ws.Range("E21").Value = "2101"
$ws.Range("E22").Value = "2012"
$ws.Range("E23").Value = "2011"

$ws.Range("F16").Value = 24578
$ws.Range("F23").Value = 35112
